$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5725536346435547
$ws.Range("B1").Value = 2.698996782302856
$ws.Range("C1").Value = 4.448913097381592
$ws.Range("D1").Value = 1.671173334121704
$ws.Range("E1").Value = 1.192523121833801
